# Fruta / hortaliza, semanal
# Insert a brand-new weekly record as row 54 (pushing the previous rows 54-77
# down to 55-78), and append one more new record as the new last row (79).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row at position 54, shifting existing rows 54:77 down to 55:78 ---
$ws.Rows.Item(54).Insert()

# Column D uses a date/time number format (style copied from the row below,
# which used to be row 54 and is now row 55, so its format is unchanged).
$dateFormat = $ws.Cells.Item(55, 4).NumberFormat()

# --- 2) Fill in the data for the newly inserted row 54 ---
$ws.Cells.Item(54, 1).Value = 6
$ws.Cells.Item(54, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(54, 3).Value = "Metropolitana"
$ws.Cells.Item(54, 4).Value = 44784
$ws.Cells.Item(54, 4).NumberFormat = $dateFormat
$ws.Cells.Item(54, 5).Value = 13
$ws.Cells.Item(54, 6).Value = "Fruta"
$ws.Cells.Item(54, 7).Value = 100108
$ws.Cells.Item(54, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(54, 9).Value = 100108007
$ws.Cells.Item(54, 10).Value = "Coco"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 100
$ws.Cells.Item(54, 14).Value = 27000
$ws.Cells.Item(54, 15).Value = 28000
$ws.Cells.Item(54, 16).Value = 27500
$ws.Cells.Item(54, 17).Value = '$/malla 20 unidades'
$ws.Cells.Item(54, 18).Value = "Perú"
$ws.Cells.Item(54, 19).Value = 1375
$ws.Cells.Item(54, 20).Value = 20

# --- 3) Append a brand-new row 79 after the (shifted) last data row 78 ---
$ws.Cells.Item(79, 1).Value = 6
$ws.Cells.Item(79, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(79, 3).Value = "Metropolitana"
$ws.Cells.Item(79, 4).Value = 44489
$ws.Cells.Item(79, 4).NumberFormat = $dateFormat
$ws.Cells.Item(79, 5).Value = 13
$ws.Cells.Item(79, 6).Value = "Fruta"
$ws.Cells.Item(79, 7).Value = 100108
$ws.Cells.Item(79, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(79, 9).Value = 100108007
$ws.Cells.Item(79, 10).Value = "Coco"
$ws.Cells.Item(79, 11).Value = "Sin especificar"
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 150
$ws.Cells.Item(79, 14).Value = 28000
$ws.Cells.Item(79, 15).Value = 30000
$ws.Cells.Item(79, 16).Value = 29000
$ws.Cells.Item(79, 17).Value = '$/malla 20 unidades'
$ws.Cells.Item(79, 18).Value = "Perú"
$ws.Cells.Item(79, 19).Value = 1450
$ws.Cells.Item(79, 20).Value = 20

Write-Host "Applied weekly Fruta/Coco update"
